{"js": "// UC10 correction: remove the \"duplicate login\" requirement row (R10-3) and\n// renumber the old \"Exit button\" requirement row (R10-4) down to R10-3.\n// Also the table's preferred width gets pinned to a fixed 9775 dxa (488.75pt)\n// instead of \"auto\" (a side effect of Word re-laying the table out on save).\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected at least one table in the document.\");\n}\nconst table = tables.items[0];\n\n// --- 1. Pin the table width to a fixed 9775 dxa (= 488.75pt), w:type=\"dxa\" ---\n// The Table.width property setter isn't wired through to the OOXML in this\n// host, so round-trip through getOoxml()/insertOoxml() to rewrite <w:tblW>.\nconst ooxmlResult = table.getOoxml();\nawait context.sync();\n\nlet tableXml = ooxmlResult.value;\nconst fixedTblW = '<w:tblW w:w=\"9775\" w:type=\"dxa\"/>';\nif (tableXml.indexOf('<w:tblW w:w=\"0\" w:type=\"auto\"/>') !== -1) {\n  tableXml = tableXml.replace('<w:tblW w:w=\"0\" w:type=\"auto\"/>', fixedTblW);\n} else {\n  // Fall back to a generic regex in case of attribute-order/whitespace drift.\n  tableXml = tableXml.replace(/<w:tblW\\s+[^>]*\\/>/, fixedTblW);\n}\ntable.insertOoxml(tableXml, Word.InsertLocation.replace);\nawait context.sync();\n\n// The table (and its rows/cells) were replaced wholesale by insertOoxml, so\n// re-resolve everything from the document body rather than reusing old\n// references.\ntables.load(\"items\");\nawait context.sync();\nconst table2 = tables.items[0];\ntable2.rows.load(\"items\");\nawait context.sync();\n\n// --- 2. Locate the old \"R10-3\" (duplicate login) row by its text content\n// (robust to any shift in row index) and delete it entirely ---\nasync function findRowByPredicate(predicate) {\n  table2.rows.load(\"items\");\n  await context.sync();\n  const currentRows = table2.rows.items;\n  for (const r of currentRows) {\n    r.load(\"values\");\n  }\n  await context.sync();\n  for (const r of currentRows) {\n    if (predicate(r.values[0][0])) {\n      return r;\n    }\n  }\n  return null;\n}\n\nconst duplicateLoginRow = await findRowByPredicate(\n  (cellText) => cellText.indexOf(\"R10-3\") === 0 && cellText.indexOf(\"already logged in\") !== -1\n);\nif (!duplicateLoginRow) {\n  throw new Error('Could not find the \"R10-3\" (duplicate login) row to remove.');\n}\nduplicateLoginRow.delete();\nawait context.sync();\n\n// --- 3. Re-resolve the old \"R10-4\" row AFTER the deletion above (deleting a\n// row shifts indices, so any reference captured beforehand is stale) and\n// relabel its heading paragraph to \"R10-3\" ---\nconst exitButtonRow = await findRowByPredicate((cellText) => cellText.indexOf(\"R10-4\") === 0);\nif (!exitButtonRow) {\n  throw new Error('Could not find the \"R10-4\" (exit button) row to renumber.');\n}\n\nexitButtonRow.cells.load(\"items\");\nawait context.sync();\nconst headingCell = exitButtonRow.cells.items[0];\nheadingCell.body.paragraphs.load(\"items\");\nawait context.sync();\nconst headingParagraph = headingCell.body.paragraphs.items[0];\nheadingParagraph.insertText(\"R10-3\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# UC10 correction: remove the \"duplicate login\" requirement row (R10-3) and\n# renumber the old \"Exit button\" requirement row (R10-4) down to R10-3.\n# Also the table's preferred width gets pinned to a fixed 9775 dxa (488.75pt)\n# instead of \"auto\" (a side effect of Word re-laying the table out on save).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# --- 1. Locate the old \"R10-3\" duplicate-login-warning row by its text\n# content (robust to any shift in row index) and delete it entirely ---\n$dupIdx = -1\nfor ($i = 1; $i -le $t.Rows.Count; $i++) {\n    $rowText = $t.Rows.Item($i).Range.Text\n    if ($rowText.StartsWith(\"R10-3\") -and $rowText.Contains(\"already logged in\")) {\n        $dupIdx = $i\n        break\n    }\n}\nif ($dupIdx -lt 0) {\n    throw \"Could not find the 'R10-3' (duplicate login) row to remove.\"\n}\n$t.Rows.Item($dupIdx).Delete()\n\n# --- 2. Re-scan for the old \"R10-4\" row AFTER the deletion above (deleting a\n# row shifts indices, so any index captured beforehand could be stale) and\n# relabel its heading paragraph to \"R10-3\" ---\n$exitIdx = -1\nfor ($i = 1; $i -le $t.Rows.Count; $i++) {\n    $rowText = $t.Rows.Item($i).Range.Text\n    if ($rowText.StartsWith(\"R10-4\")) {\n        $exitIdx = $i\n        break\n    }\n}\nif ($exitIdx -lt 0) {\n    throw \"Could not find the 'R10-4' (exit button) row to renumber.\"\n}\n$exitRow = $t.Rows.Item($exitIdx)\n$headingPara = $exitRow.Cells.Item(1).Range.Paragraphs.Item(1)\n$headingPara.Range.Text = \"R10-3\"\n\n# --- 3. Pin the table width to a fixed 9775 dxa (= 488.75pt) ---\n$t.PreferredWidthType = 3\n$t.PreferredWidth = 488.75\n"}
